# Update "杭州-漫展信息.xlsx" to the newer scraped snapshot.
#
# Sheet "展览" (展览 = Exhibitions, the 1st / active sheet):
#   * A new event ("杭州·逐月节·园游会...") happened on 2024-10-04, which sorts
#     before the existing 2024-10-05 "华彩的摔跤宴" row, so it is inserted as a
#     new row 30 and everything from the old row 30 onward shifts down by one.
#   * Several "want to go" / price counters (F/G columns) were refreshed to
#     newer scraped values, both in rows unaffected by the insertion (rows
#     1-29) and in the rows that got shifted down (new rows 31-48).
#
# Sheet "全部类型" (全部类型 = All types, the 4th sheet) has no row insertion
# (the new event doesn't appear there before or after), just the same kind of
# refreshed F/G counters, in place.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# --- F-column refresh for rows not affected by the insertion (rows 1-29) ---
$ws1.Range("F3").Value = 8853
$ws1.Range("F7").Value = 2115
$ws1.Range("F8").Value = 590
$ws1.Range("F14").Value = 75
$ws1.Range("F16").Value = 8740
$ws1.Range("F17").Value = 165
$ws1.Range("F19").Value = 201
$ws1.Range("F20").Value = 116
$ws1.Range("F21").Value = 1826
$ws1.Range("F28").Value = 197

# --- Insert the new row 30 (shifts old rows 30-47 down to 31-48) ---
$ws1.Rows.Item(30).Insert()

# Column A just holds "row number - 1" as a literal value; row-insert leaves
# the new row blank, so give it the right label + copy the header style of
# the surrounding A-column cells (bold / centered / bordered).
$ws1.Range("A31").Copy($ws1.Range("A30"))
$ws1.Range("A30").Value = 29

$ws1.Range("B30").NumberFormat = "@"
$ws1.Range("B30").Value = "2024-10-04"
$ws1.Range("C30").Value = "杭州·逐月节·园游会·原神×绝区零×崩铁×崩坏同人only"
$ws1.Range("D30").Value = "莫干山路987号 资辉壹方汇"
$ws1.Range("E30").Value = "2024.10.04 09:30-10.05 17:00"
$ws1.Range("F30").Value = 1
$ws1.Range("G30").Value = 58
$ws1.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=92406"
$ws1.Range("I30").Value = "//i1.hdslb.com/bfs/openplatform/202409/mQh43oPd1726134932363.png"

# --- F/G refresh for the rows that shifted down (now rows 31-48) ---
$ws1.Range("F31").Value = 31
$ws1.Range("F32").Value = 39
$ws1.Range("F33").Value = 26
$ws1.Range("F34").Value = 27
$ws1.Range("F35").Value = 2212
$ws1.Range("G35").Value = 49.9
$ws1.Range("F36").Value = 869
$ws1.Range("F37").Value = 514
$ws1.Range("F38").Value = 5
$ws1.Range("F39").Value = 6
$ws1.Range("F40").Value = 10
$ws1.Range("F41").Value = 243
$ws1.Range("F42").Value = 175
$ws1.Range("F43").Value = 7
$ws1.Range("F44").Value = 735
$ws1.Range("F45").Value = 81
$ws1.Range("F46").Value = 98
$ws1.Range("F47").Value = 69
$ws1.Range("F48").Value = 3987

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (no row insertion here, just refreshed counters)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value = 8853
$ws4.Range("F9").Value = 2115
$ws4.Range("F10").Value = 590
$ws4.Range("F17").Value = 75
$ws4.Range("F19").Value = 8740
$ws4.Range("F20").Value = 165
$ws4.Range("F22").Value = 201
$ws4.Range("F23").Value = 1826
$ws4.Range("F28").Value = 197
$ws4.Range("F31").Value = 39
$ws4.Range("F33").Value = 27
$ws4.Range("F34").Value = 2212
$ws4.Range("G34").Value = 49.9
$ws4.Range("F37").Value = 514
$ws4.Range("F39").Value = 243
$ws4.Range("F41").Value = 175
$ws4.Range("F42").Value = 69
